$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.871.46"
$ws.Range("E2").Value = "  -1.00%  "

$ws.Range("D3").Value = "1.737.12"
$ws.Range("E3").Value = "  -1.89%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.30"
$ws.Range("E5").Value = "  -3.72%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9998"
$ws.Range("E6").Value = "  -0.02%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5257"
$ws.Range("E7").Value = "  +0.39%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2767"
$ws.Range("E8").Value = "  -0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.37"
$ws.Range("E9").Value = "  -2.91%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06143"
$ws.Range("E10").Value = "  -1.28%  "

$ws.Range("D11").Value = "1.738.63"
$ws.Range("E11").Value = "  -1.84%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07086"
$ws.Range("E12").Value = "  +0.74%  "

$ws.Range("E13").Value = "  -6.24%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6424"
$ws.Range("E14").Value = "  -1.12%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.529"
$ws.Range("E15").Value = "  -0.20%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.98"
$ws.Range("E16").Value = "  -1.75%  "

$ws.Range("E17").Value = "  +0.00%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9998"
$ws.Range("E18").Value = "  -0.02%  "

$ws.Range("D19").Value = "25.859.04"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.54"
$ws.Range("E20").Value = "  -1.99%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006650"
$ws.Range("E21").Value = "  -1.98%  "

$ws.Range("D22").Value = "1.959.61"
$ws.Range("E22").Value = "  -2.00%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.253"
$ws.Range("E23").Value = "  +3.96%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.820"
$ws.Range("E24").Value = "  +4.28%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.164"
$ws.Range("E25").Value = "  -1.00%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.22"
$ws.Range("E26").Value = "  +1.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.512"
$ws.Range("E27").Value = "  +1.40%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.25"
$ws.Range("E28").Value = "  +0.29%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.785"
$ws.Range("E29").Value = "  -5.06%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.50"
$ws.Range("E30").Value = "  -0.68%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08328"
$ws.Range("E31").Value = "  -1.19%  "

$ws.Range("E32").Value = "  -0.76%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.536"
$ws.Range("E33").Value = "  +1.63%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04510"
$ws.Range("E34").Value = "  +0.73%  "

$ws.Range("E35").Value = "  -1.81%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9733"
$ws.Range("E36").Value = "  -3.74%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6193"
$ws.Range("E37").Value = "  +1.48%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.679"
$ws.Range("E38").Value = "  -3.26%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01576"
$ws.Range("E39").Value = "  -0.86%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.912"
$ws.Range("E40").Value = "  -4.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9994"
$ws.Range("E41").Value = "  -0.22%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.98"
$ws.Range("E42").Value = "  -3.07%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3855"
$ws.Range("E43").Value = "  -0.74%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.022"
$ws.Range("E44").Value = "  +1.51%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.7241"
$ws.Range("E45").Value = "  -3.74%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05324"
$ws.Range("E46").Value = "  -3.45%  "

$ws.Range("E47").Value = "  -0.08%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.208"
$ws.Range("E48").Value = "  -3.38%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.43"
$ws.Range("E49").Value = "  +0.87%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "30.03"
$ws.Range("E50").Value = "  -1.13%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.629"
$ws.Range("E51").Value = "  +2.10%  "
